$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92, shifting existing rows 92..161 down to 93..162.
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with its values.
$ws.Range("A92").Value = 4
$ws.Range("B92").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C92").Value = "Los Lagos"
$ws.Range("D92").Value = 44574
$ws.Range("E92").Value = 10
$ws.Range("F92").Value = 100112039
$ws.Range("G92").Value = "Ciboulette"
$ws.Range("H92").Value = "Sin especificar"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 120
$ws.Range("K92").Value = 3000
$ws.Range("L92").Value = 3000
$ws.Range("M92").Value = 3000
$ws.Range("N92").Value = "$/docena de atados"
$ws.Range("O92").Value = "Región Metropolitana"
$ws.Range("P92").Value = 1000
$ws.Range("Q92").Value = 3
$ws.Range("R92").Value = "Hortaliza"

# Make sure the worksheet dimension covers the new row.
$ws.Range("A1:R162").Select()
